$ws = $excel.ActiveWorkbook.ActiveSheet

# Row, Column, NewValue triples (Row=spreadsheet row, Column=spreadsheet column index)
$updates = @"
2,2,1.02
2,3,1.035383042756525
2,4,1.038027399199306
2,5,1.043239394799286
2,6,1.051323590264847
2,9,1.03248832586516
2,10,1.040497118508461
2,11,1.040816602684436
2,12,1.046013820906619
2,13,1.05407540986662
2,14,1.041974743763431
3,2,1.02
3,3,1.03644764337012
3,4,1.039006121426829
3,5,1.04422181502091
3,6,1.052474633873844
3,9,1.032619345340919
3,10,1.041204663579154
3,11,1.041604807522457
3,12,1.046806787332633
3,13,1.055038206446775
3,14,1.04268329362928
4,2,1.02
4,3,1.037136584836365
4,4,1.039639772019352
4,5,1.044857983438921
4,6,1.053220228044359
4,9,1.032702568322428
4,10,1.041662004229239
4,11,1.04211456682116
4,12,1.047319751683657
4,13,1.055661410622551
4,14,1.043141283755546
5,2,1.02
5,3,1.037426233672469
5,4,1.039906242637744
5,5,1.045125542503994
5,6,1.053533865390889
5,9,1.032737182612707
5,10,1.041854153111253
5,11,1.042328806666146
5,12,1.047535368662868
5,13,1.055923455913031
5,14,1.043333705511013
6,2,1.02
6,3,1.037474868038942
6,4,1.039950989138061
6,5,1.045170473516724
6,6,1.053586537595899
6,9,1.032742972643478
6,10,1.041886408879203
6,11,1.042364774779145
6,12,1.047571569747551
6,13,1.055967457413976
6,14,1.043366007085853
7,2,1.02
7,3,1.037140455067918
7,4,1.039643332283506
7,5,1.044861558129216
7,6,1.05322441813716
7,9,1.032703032304896
7,10,1.041664572192519
7,11,1.042117429752187
7,12,1.047322632899139
7,13,1.055664911887229
7,14,1.043143855365628
8,2,1.02
8,3,1.035742814384576
8,4,1.03835808990527
8,5,1.043571309402474
8,6,1.051712426827406
8,9,1.032532926388134
8,10,1.040736337655211
8,11,1.041083034282093
8,12,1.046281836061501
8,13,1.054400747841787
8,14,1.042214302628808
9,2,1.02
9,3,1.033280547668975
9,4,1.03609603753323
9,5,1.041301390920452
9,6,1.049054172190726
9,9,1.032221278319526
9,10,1.039096939226941
9,11,1.039258302055975
9,12,1.044446771311925
9,13,1.05217474547819
9,14,1.040572576066728
10,2,1.02
10,3,1.031639393141106
10,4,1.03458984025308
10,5,1.039790590211283
10,6,1.047286080647826
10,9,1.032005528532655
10,10,1.038001509158192
10,11,1.038040484739569
10,12,1.04322269989884
10,13,1.050691834158889
10,14,1.039475590361634
11,2,1.02
11,3,1.030928835114776
11,4,1.033938078659503
11,5,1.039136986839061
11,6,1.046521443543061
11,9,1.031910215866671
11,10,1.037526584878365
11,11,1.037512842969327
11,12,1.042692499887098
11,13,1.050049975640871
11,14,1.038999991634856
12,2,1.02
12,3,1.030664912644016
12,4,1.033696050362646
12,5,1.038894297160611
12,6,1.046237567018209
12,9,1.031874528516439
12,10,1.037350087363614
12,11,1.037316805473466
12,12,1.042495534610114
12,13,1.049811598529798
12,14,1.038823243473389
13,2,1.02
13,3,1.030721524446586
13,4,1.033747963311717
13,5,1.038946350949854
13,6,1.046298452995311
13,9,1.031882196421686
13,10,1.037387950728919
13,11,1.03735885834623
13,12,1.042537785479278
13,13,1.04986272952425
13,14,1.038861160609015
14,2,1.02
14,3,1.030907018991412
14,4,1.0339180711988
14,5,1.039116924244692
14,6,1.046497975289933
14,9,1.031907271731124
14,10,1.037511997364183
14,11,1.03749663942457
14,12,1.042676219193674
14,13,1.050030270545689
14,14,1.038985383404729
15,2,1.02
15,3,1.031021309667921
15,4,1.034022888895907
15,5,1.039222031706869
15,6,1.046620926601363
15,9,1.031922683832635
15,10,1.037588414732538
15,11,1.037581524555374
15,12,1.04276150941039
15,13,1.050133503093442
15,14,1.039061909294514
16,2,1.02
16,3,1.031686552025853
16,4,1.034633104629811
16,5,1.039833980010957
16,6,1.047336847263784
16,9,1.032011814284548
16,10,1.038033015767233
16,11,1.038075495873081
16,12,1.043257883967121
16,13,1.05073443749121
16,14,1.039507141713674
17,2,1.02
17,3,1.032103860091134
17,4,1.035015992512892
17,5,1.040217995239114
17,6,1.047786181764753
17,9,1.032067217177574
17,10,1.038311742747528
17,11,1.038385265705957
17,12,1.043569201368033
17,13,1.0511114551107
17,14,1.039786264518251
18,2,1.02
18,3,1.032347275875067
18,4,1.035239366189139
18,5,1.040442041099713
18,6,1.048048363444043
18,9,1.032099350223401
18,10,1.03847426198865
18,11,1.038565918442053
18,12,1.043750771233226
18,13,1.051331387419163
18,14,1.039949014555337
19,2,1.02
19,3,1.032430275578162
19,4,1.035315537901779
19,5,1.040518444491074
19,6,1.048137776311946
19,9,1.032110275796042
19,10,1.03852966709778
19,11,1.038627511136443
19,12,1.043812679098804
19,13,1.051406382779766
19,14,1.040004498346077
20,2,1.02
20,3,1.032059086159413
20,4,1.034974907957479
20,5,1.040176788217929
20,6,1.047737962862395
20,9,1.032061291849301
20,10,1.038281843926644
20,11,1.038352033516365
20,12,1.043535801643821
20,13,1.051071002174466
20,14,1.039756323237612
21,2,1.02
21,3,1.030852395187768
21,4,1.033867976871269
21,5,1.03906669223876
21,6,1.046439216980905
21,9,1.031899895515122
21,10,1.037475471192069
21,11,1.037456067651554
21,12,1.042635454615863
21,13,1.049980932850387
21,14,1.038948805361262
22,2,1.02
22,3,1.030093760705561
22,4,1.033172381364427
22,5,1.038369238372976
22,6,1.045623475434673
22,9,1.031796776302293
22,10,1.036967954893066
22,11,1.036892461536976
22,12,1.04206922383575
22,13,1.049295781187078
22,14,1.038440568330901
23,2,1.02
23,3,1.030495921677437
23,4,1.033541094026458
23,5,1.038738923728054
23,6,1.046055836714901
23,9,1.031851597410051
23,10,1.037237047954488
23,11,1.037191266023562
23,12,1.042369407404975
23,13,1.049658972405424
23,14,1.038710043535336
24,2,1.02
24,3,1.032079317549858
24,4,1.034993472170786
24,5,1.040195407723063
24,6,1.047759750622913
24,9,1.032063969814047
24,10,1.038295354094338
24,11,1.038367049807321
24,12,1.04355089359178
24,13,1.051089281039797
24,14,1.039769852591294
25,2,1.02
25,3,1.033917038585475
25,4,1.036680509587486
25,5,1.041887783495319
25,6,1.049740675845434
25,9,1.032303255311862
25,10,1.0395212040298
25,11,1.039730274248482
25,12,1.044921302715074
25,13,1.05407540986662
25,14,1.040997443374271
"@

$lines = $updates -split "`n" | Where-Object { $_.Trim() -ne "" }
foreach ($line in $lines) {
    $parts = $line.Trim() -split ","
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
}
